# Auto-generated edit script: update F (想去人数) and G (最低票价) values
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 338
$ws.Range("F3").Value = 3436
$ws.Range("F5").Value = 8042
$ws.Range("F7").Value = 49
$ws.Range("F8").Value = 2032
$ws.Range("F9").Value = 1897
$ws.Range("F10").Value = 45
$ws.Range("F12").Value = 1074
$ws.Range("G13").Value = 58.8
$ws.Range("F14").Value = 52
$ws.Range("F16").Value = 21
$ws.Range("F19").Value = 140
$ws.Range("F20").Value = 1096
$ws.Range("F21").Value = 683
$ws.Range("F23").Value = 42
$ws.Range("F26").Value = 4027
$ws.Range("F27").Value = 48
$ws.Range("F28").Value = 43494
$ws.Range("F29").Value = 3788
$ws.Range("F31").Value = 969
$ws.Range("F32").Value = 641
$ws.Range("F34").Value = 802
$ws.Range("F36").Value = 550
$ws.Range("F37").Value = 167
$ws.Range("F39").Value = 550
$ws.Range("F41").Value = 834
$ws.Range("F46").Value = 61
$ws.Range("F47").Value = 28
$ws.Range("F48").Value = 9
$ws.Range("F49").Value = 2434

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 223
$ws.Range("F5").Value = 47
$ws.Range("F6").Value = 47
$ws.Range("F7").Value = 41
$ws.Range("F11").Value = 97
$ws.Range("F14").Value = 29
$ws.Range("F15").Value = 69
$ws.Range("F18").Value = 137
$ws.Range("F19").Value = 7175
$ws.Range("F25").Value = 15
$ws.Range("F27").Value = 91
$ws.Range("F35").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2160
$ws.Range("F5").Value = 1450
$ws.Range("F8").Value = 2294
$ws.Range("F9").Value = 9173
$ws.Range("F10").Value = 1440

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 338
$ws.Range("F3").Value = 2160
$ws.Range("F5").Value = 1450
$ws.Range("F6").Value = 1440
$ws.Range("F9").Value = 49
$ws.Range("F10").Value = 45
$ws.Range("F11").Value = 1074
$ws.Range("G12").Value = 58.8
$ws.Range("F13").Value = 52
$ws.Range("F15").Value = 21
$ws.Range("F17").Value = 223
$ws.Range("F18").Value = 47
$ws.Range("F20").Value = 140
$ws.Range("F21").Value = 1096
$ws.Range("F22").Value = 683
$ws.Range("F23").Value = 42
$ws.Range("F25").Value = 4027
$ws.Range("F26").Value = 48
$ws.Range("F27").Value = 97
$ws.Range("F30").Value = 29
$ws.Range("F31").Value = 969
$ws.Range("F33").Value = 802
$ws.Range("F34").Value = 550
$ws.Range("F35").Value = 69
$ws.Range("F36").Value = 167
$ws.Range("F37").Value = 550
$ws.Range("F39").Value = 834
$ws.Range("F43").Value = 62
$ws.Range("F44").Value = 91
$ws.Range("F46").Value = 28
$ws.Range("F48").Value = 2434
$ws.Range("F50").Value = 11
